$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A112").Value = 'female_int_eligibility'
$ws.Range("A113").Value = 'year_of birth'
$ws.Range("A114").Value = 'individual_file_pregnancy_status'
$ws.Range("A115").Value = 'brand_of_net'
$ws.Range("A116").Value = 'century_day_code_of_birth'
$ws.Range("A117").Value = 'month_of_data_collection'
$ws.Range("A118").Value = 'hemoglobin_level_adjusted_for_altitude_g_dl'
$ws.Range("A119").Value = 'na_read_consent_statement_for_malaria'
$ws.Range("A120").Value = 'insecticide_treated_net'
$ws.Range("A121").Value = 'childs_age_in_days'
$ws.Range("A122").Value = 'read_consent_statement_hemoglobin'
$ws.Range("A123").Value = 'index_to_household_schedule_hmhidx'
$ws.Range("A124").Value = 'month_of_birth'
$ws.Range("A125").Value = 'mothers_line_number'
$ws.Range("A126").Value = 'malariae_present'
$ws.Range("A127").Value = 'net_observed_by_interviewer'
$ws.Range("A128").Value = 'rshp_to_head'
$ws.Range("A129").Value = 'date_measured_day'
$ws.Range("A130").Value = 'line_number_of_person_slept_in_net'
$ws.Range("A131").Value = 'net_from_antenatal_immunization_visit'
$ws.Range("A132").Value = 'childs_age_in_months'
$ws.Range("A133").Value = 'sex_of_member'
$ws.Range("A134").Value = 'line_number_of_parent_caretaker'
$ws.Range("A135").Value = 'net_design_no'
$ws.Range("A136").Value = 'slept_under_net'
$ws.Range("A137").Value = 'corr_age'
$ws.Range("A138").Value = 'flag_age'
$ws.Range("A139").Value = 'final_blood_smear_test'
$ws.Range("A140").Value = 'date_measured_month'
$ws.Range("A141").Value = 'date_measured_year'
$ws.Range("A142").Value = 'caretaker_line_number'
$ws.Range("A143").Value = 'vivax_present'
$ws.Range("A144").Value = 'ovale_present'
$ws.Range("A145").Value = 'sex'
$ws.Range("A146").Value = 'age_of_member'
$ws.Range("A147").Value = 'anemia_level'
$ws.Range("A148").Value = 'malaria_measurement_result'
$ws.Range("A149").Value = 'day_of_birth'
$ws.Range("A150").Value = 'day_of_data_collection'
$ws.Range("A151").Value = 'line_number'
$ws.Range("A152").Value = 'childs_age_in_months_country_specific'
$ws.Range("A153").Value = 'number_of_persons_slept_under_net'
$ws.Range("A154").Value = 'cmc_date_of_birth'
$ws.Range("A155").Value = 'months_ago_net_obtained'
$ws.Range("A156").Value = 'hemoglobin_level_g_dl'
$ws.Range("A157").Value = 'child_age_in_months'
$ws.Range("A158").Value = 'century_day_code_of_measurement'
$ws.Range("A159").Value = 'fieldworker_measurer_code'
$ws.Range("A160").Value = 'mosquito_bed_net_designation_number'
$ws.Range("A161").Value = 'blood_smear_bar_code'
$ws.Range("A162").Value = 'bed_net_type'
$ws.Range("A163").Value = 'index_to_household_schedule_hc0'
$ws.Range("A164").Value = 'childs_age_in_days_country_specific'
$ws.Range("A165").Value = 'someone_slept_under_net_last_night'
$ws.Range("A166").Value = 'childs_age_in_months_country_specific_hml16a'
$ws.Range("A167").Value = 'year_of_data_collection'
$ws.Range("A168").Value = 'children_hemoglobin_elig'
$ws.Range("A169").Value = 'falciparum_present'
$ws.Range("A170").Value = 'result_of_measurement_hemoglobin'
$ws.Range("A171").Value = 'usual_resident'
$ws.Range("A172").Value = 'slept_last_night'
$ws.Range("A173").Value = 'slept_llin_net'
$ws.Range("A174").Value = 'completeleness_of_hc32_info'
$ws.Range("A175").Value = 'rapid_test_result'
$ws.Range("A176").Value = 'fieldworker_malaria_measurer_code'
